$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C32 (en_US): curly double quotes -> straight single quotes around "Diαbolic Crisis?"
$ws.Range("C32").Value = "[name=""Muelsyse""]The 'Diαbolic Crisis?'`n"

# C34 (en_US): curly double quotes -> straight single quotes around "Diαbolic Crisis"
$ws.Range("C34").Value = "[name=""Silence""]The 'Diαbolic Crisis' precipitated Saria's departure. What happened after that, I presume you understand even better than I do, Director Muelsyse.`n"

# D54 (ko_KR): Muelsyse's first "Huh?" changes from 네? to 에?
$ws.Range("D54").Value = "[name=""뮤엘시스""]에?`n"

# D56 (ko_KR): Muelsyse's second "Huh?" (previously duplicated the same shared
# string as D54) now gets its own, distinct, more emphatic line.
$ws.Range("D56").Value = "[name=""뮤엘시스""]에?!`n"
